$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Fix header label: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Remove the two empty "section header" rows (these previously held
# "situação do domicílio" and "grandes regiões e unidades da federação"
# labels with no accompanying data); data below shifts up to fill them.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
